$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D. Excel shifts the old D:K data to E:L,
# but Insert() only picks up the left-neighbour's format for the brand
# new column, so we explicitly re-stamp formats for D from the (just
# shifted) E column - which still carries the original D formatting.
$ws.Columns("D").Insert()
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New FY2018 ("12/31/2018") column of data, inserted ahead of the
# existing years.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1038200
$ws.Range("D9").Value = 306800
$ws.Range("D10").Value = 731400
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 3200
$ws.Range("D15").Value = 325400
$ws.Range("D17").Value = 681200
$ws.Range("D18").Value = 357000
$ws.Range("D20").Value = 67400
$ws.Range("D21").Value = 749800
$ws.Range("D22").Value = 263700
$ws.Range("D23").Value = 160700
$ws.Range("D24").Value = 86700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 74000
$ws.Range("D27").Value = 56100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -67400
$ws.Range("D33").Value = 56100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 56100
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 841600
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 157300
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 110600
$ws.Range("D46").Value = 1109500
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 5005600
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 70300
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 6185400
$ws.Range("D57").Value = 25700
$ws.Range("D58").Value = 162900
$ws.Range("D59").Value = 244000
$ws.Range("D60").Value = 432600
$ws.Range("D61").Value = 2896200
$ws.Range("D62").Value = 142400
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 4855700
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1329700
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 56100
$ws.Range("D83").Value = 325400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 434100
$ws.Range("D91").Value = -23400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -23400
$ws.Range("D96").Value = -55400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -416700
$ws.Range("D101").Value = -1000
$ws.Range("D102").Value = -7000

# Row 14 ("Non Recurring") lost its historical 0-values for FY2011-2016
# in this refresh - those columns now read "NA" and only the oldest
# (FY2011->now column K) and newest (FY2018->now column D) columns keep
# numbers.
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "NA"
$ws.Range("H14").Value = "NA"
$ws.Range("I14").Value = "NA"
$ws.Range("J14").Value = "NA"
